$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the values in range B2:D9 to 0, matching the diff which zeros out
# all the correlation values in that range.
$ws.Range("B2:D9").Value = 0
